# Edit script: insert a new "15_02_01_05_*" (three-wheeler / mini vehicle?) fuel-type
# breakdown block into the sector reference table, mirroring the pattern already used
# for the other "15_02_0x_0y" engine-type groups.
#
# Net effect vs. the original workbook:
#   - 7 new rows are inserted right before the old row 197 ("15_02_02_freight"),
#     pushing the existing "15_02_02_freight" ... "15_02_02_04_08_fuel_cell_ev" block
#     down by 7 rows (old rows 197-280 become new rows 204-287).
#   - 1 new row is appended at the very end of the table (new row 288) to hold the
#     matching "..._08_fuel_cell_ev" entry for this new group.
#   - The new rows follow the existing convention: column B and C both hold the same
#     sector code text, column D holds "(new)" (except for the very last appended
#     row, which only has B and C populated, matching the rest of the "_08_fuel_cell_ev"
#     rows that also lack a D value).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 7 blank rows at row 197, shifting everything below it down.
$ws.Rows("197:203").Insert()

# Append the new "...05_08_fuel_cell_ev" row at the bottom FIRST, so that its shared
# string is registered before the rows below, matching workbook authoring order.
$ws.Range("B288").Value = "15_02_01_05_08_fuel_cell_ev"
$ws.Range("C288").Value = "15_02_01_05_08_fuel_cell_ev"

# Fill in the newly inserted rows 197-203 with the new sub-group codes.
$newCodes = @(
    "15_02_01_05_01_diesel_engine",
    "15_02_01_05_02_gasoline_engine",
    "15_02_01_05_03_battery_ev",
    "15_02_01_05_04_compressed_natual_gas",
    "15_02_01_05_05_plugin_hybrid_ev_gasoline",
    "15_02_01_05_06_plugin_hybrid_ev_diesel",
    "15_02_01_05_07_liquified_petroleum_gas"
)

for ($i = 0; $i -lt $newCodes.Length; $i++) {
    $row = 197 + $i
    $code = $newCodes[$i]
    $ws.Range("B$row").Value = $code
    $ws.Range("C$row").Value = $code
    $ws.Range("D$row").Value = "(new)"
}

# Put the active selection where the original author left it after the edit.
$ws.Range("B285").Select()
